$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-28 09:51:12"
$wsZhCn.Range("D3").Value = "2016-01-28 09:51:12"
$wsZhCn.Range("G2").Value = "2016-01-28 09:52:02"
$wsZhCn.Range("G3").Value = "2016-01-28 09:52:02"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-28 09:51:25"
$wsDeDe.Range("D3").Value = "2016-01-28 09:51:25"
$wsDeDe.Range("G2").Value = "2016-01-28 09:52:25"
$wsDeDe.Range("G3").Value = "2016-01-28 09:52:25"
